$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This worksheet lists weekly Haba (fava bean) price records for the
# "Vega Monumental Concepción" market. This update refreshes the
# weekly data snapshot: the Fecha (date), Volumen, Precio mínimo/máximo,
# Precio promedio ponderado, Origen and Precio $/Kg columns are updated
# per row to reflect the latest weekly price report.

    # Row 2 <- values from original row 28
    $ws.Range("D2").Value = 44483
    $ws.Range("J2").Value = 350
    $ws.Range("K2").Value = 5500
    $ws.Range("L2").Value = 6000
    $ws.Range("M2").Value = 5714
    $ws.Range("O2").Value = 'Región Metropolitana'
    $ws.Range("P2").Value = 229
    # Row 3 <- values from original row 20
    $ws.Range("D3").Value = 44351
    $ws.Range("J3").Value = 100
    $ws.Range("K3").Value = 15000
    $ws.Range("L3").Value = 16000
    $ws.Range("M3").Value = 15500
    $ws.Range("O3").Value = 'Región Metropolitana'
    $ws.Range("P3").Value = 620
    # Row 4 <- values from original row 21
    $ws.Range("D4").Value = 44509
    $ws.Range("J4").Value = 100
    $ws.Range("K4").Value = 6500
    $ws.Range("L4").Value = 7000
    $ws.Range("M4").Value = 6750
    $ws.Range("O4").Value = 'Región Metropolitana'
    $ws.Range("P4").Value = 270
    # Row 5 <- values from original row 23
    $ws.Range("D5").Value = 44673
    $ws.Range("J5").Value = 80
    $ws.Range("K5").Value = 18000
    $ws.Range("L5").Value = 19000
    $ws.Range("M5").Value = 18375
    $ws.Range("O5").Value = 'Región Metropolitana'
    $ws.Range("P5").Value = 735
    # Row 6 <- values from original row 27
    $ws.Range("D6").Value = 44162
    $ws.Range("J6").Value = 100
    $ws.Range("K6").Value = 7500
    $ws.Range("L6").Value = 8000
    $ws.Range("M6").Value = 7750
    $ws.Range("O6").Value = 'Región Metropolitana'
    $ws.Range("P6").Value = 310
    # Row 7 <- values from original row 11
    $ws.Range("D7").Value = 44188
    $ws.Range("J7").Value = 100
    $ws.Range("K7").Value = 18000
    $ws.Range("L7").Value = 20000
    $ws.Range("M7").Value = 19000
    $ws.Range("O7").Value = 'Región Metropolitana'
    $ws.Range("P7").Value = 760
    # Row 8 <- values from original row 9
    $ws.Range("D8").Value = 44526
    $ws.Range("J8").Value = 100
    $ws.Range("K8").Value = 7500
    $ws.Range("L8").Value = 8000
    $ws.Range("M8").Value = 7750
    $ws.Range("O8").Value = 'Región Metropolitana'
    $ws.Range("P8").Value = 310
    # Row 9 <- values from original row 29
    $ws.Range("D9").Value = 44769
    $ws.Range("J9").Value = 100
    $ws.Range("K9").Value = 18000
    $ws.Range("L9").Value = 20000
    $ws.Range("M9").Value = 19000
    $ws.Range("O9").Value = 'Región de Coquimbo'
    $ws.Range("P9").Value = 760
    # Row 10 <- values from original row 5
    $ws.Range("D10").Value = 44537
    $ws.Range("J10").Value = 160
    $ws.Range("K10").Value = 8500
    $ws.Range("L10").Value = 9000
    $ws.Range("M10").Value = 8719
    $ws.Range("O10").Value = 'Región del Maule'
    $ws.Range("P10").Value = 349
    # Row 11 <- values from original row 16
    $ws.Range("D11").Value = 44505
    $ws.Range("J11").Value = 180
    $ws.Range("K11").Value = 6000
    $ws.Range("L11").Value = 6500
    $ws.Range("M11").Value = 6222
    $ws.Range("O11").Value = 'Región del Maule'
    $ws.Range("P11").Value = 249
    # Row 12 <- values from original row 15
    $ws.Range("D12").Value = 44467
    $ws.Range("J12").Value = 100
    $ws.Range("K12").Value = 8000
    $ws.Range("L12").Value = 9000
    $ws.Range("M12").Value = 8500
    $ws.Range("O12").Value = 'Región Metropolitana'
    $ws.Range("P12").Value = 340
    # Row 13 <- values from original row 10
    $ws.Range("D13").Value = 44482
    $ws.Range("J13").Value = 430
    $ws.Range("K13").Value = 8000
    $ws.Range("L13").Value = 8500
    $ws.Range("M13").Value = 8267
    $ws.Range("O13").Value = 'Región de O''Higgins'
    $ws.Range("P13").Value = 331
    # Row 14 <- values from original row 13
    $ws.Range("D14").Value = 44755
    $ws.Range("J14").Value = 100
    $ws.Range("K14").Value = 16000
    $ws.Range("L14").Value = 17000
    $ws.Range("M14").Value = 16500
    $ws.Range("O14").Value = 'Región de Coquimbo'
    $ws.Range("P14").Value = 660
    # Row 15 <- values from original row 14
    $ws.Range("D15").Value = 44160
    $ws.Range("J15").Value = 100
    $ws.Range("K15").Value = 9000
    $ws.Range("L15").Value = 10000
    $ws.Range("M15").Value = 9500
    $ws.Range("O15").Value = 'Región Metropolitana'
    $ws.Range("P15").Value = 380
    # Row 16 <- values from original row 31
    $ws.Range("D16").Value = 44517
    $ws.Range("J16").Value = 130
    $ws.Range("K16").Value = 6000
    $ws.Range("L16").Value = 6500
    $ws.Range("M16").Value = 6269
    $ws.Range("O16").Value = 'Región Metropolitana'
    $ws.Range("P16").Value = 251
    # Row 17 <- values from original row 7
    $ws.Range("D17").Value = 44461
    $ws.Range("J17").Value = 100
    $ws.Range("K17").Value = 13000
    $ws.Range("L17").Value = 14000
    $ws.Range("M17").Value = 13500
    $ws.Range("O17").Value = 'Provincia del Elquí'
    $ws.Range("P17").Value = 540
    # Row 18 <- values from original row 3
    $ws.Range("D18").Value = 44540
    $ws.Range("J18").Value = 140
    $ws.Range("K18").Value = 11000
    $ws.Range("L18").Value = 12000
    $ws.Range("M18").Value = 11429
    $ws.Range("O18").Value = 'Región del Maule'
    $ws.Range("P18").Value = 457
    # Row 19 <- values from original row 24
    $ws.Range("D19").Value = 44503
    $ws.Range("J19").Value = 250
    $ws.Range("K19").Value = 9000
    $ws.Range("L19").Value = 10000
    $ws.Range("M19").Value = 9400
    $ws.Range("O19").Value = 'Provincia de Melipilla'
    $ws.Range("P19").Value = 376
    # Row 20 <- values from original row 17
    $ws.Range("D20").Value = 44498
    $ws.Range("J20").Value = 220
    $ws.Range("K20").Value = 7000
    $ws.Range("L20").Value = 7500
    $ws.Range("M20").Value = 7273
    $ws.Range("O20").Value = 'Región Metropolitana'
    $ws.Range("P20").Value = 291
    # Row 21 <- values from original row 6
    $ws.Range("D21").Value = 44545
    $ws.Range("J21").Value = 140
    $ws.Range("K21").Value = 14000
    $ws.Range("L21").Value = 15000
    $ws.Range("M21").Value = 14429
    $ws.Range("O21").Value = 'Provincia de Chacabuco'
    $ws.Range("P21").Value = 577
    # Row 22 <- values from original row 25
    $ws.Range("D22").Value = 44384
    $ws.Range("J22").Value = 100
    $ws.Range("K22").Value = 12000
    $ws.Range("L22").Value = 13000
    $ws.Range("M22").Value = 12500
    $ws.Range("O22").Value = 'Región de Coquimbo'
    $ws.Range("P22").Value = 500
    # Row 23 <- values from original row 18
    $ws.Range("D23").Value = 44692
    $ws.Range("J23").Value = 100
    $ws.Range("K23").Value = 20000
    $ws.Range("L23").Value = 22000
    $ws.Range("M23").Value = 21000
    $ws.Range("O23").Value = 'Región Metropolitana'
    $ws.Range("P23").Value = 840
    # Row 24 <- values from original row 4
    $ws.Range("D24").Value = 44335
    $ws.Range("J24").Value = 100
    $ws.Range("K24").Value = 18000
    $ws.Range("L24").Value = 20000
    $ws.Range("M24").Value = 19000
    $ws.Range("O24").Value = 'Provincia de Limarí'
    $ws.Range("P24").Value = 760
    # Row 25 <- values from original row 8
    $ws.Range("D25").Value = 44316
    $ws.Range("J25").Value = 100
    $ws.Range("K25").Value = 16000
    $ws.Range("L25").Value = 18000
    $ws.Range("M25").Value = 17000
    $ws.Range("O25").Value = 'Región Metropolitana'
    $ws.Range("P25").Value = 680
    # Row 26 <- values from original row 12
    $ws.Range("D26").Value = 44523
    $ws.Range("J26").Value = 100
    $ws.Range("K26").Value = 9000
    $ws.Range("L26").Value = 10000
    $ws.Range("M26").Value = 9500
    $ws.Range("O26").Value = 'Región Metropolitana'
    $ws.Range("P26").Value = 380
    # Row 27 <- values from original row 2
    $ws.Range("D27").Value = 44476
    $ws.Range("J27").Value = 100
    $ws.Range("K27").Value = 7000
    $ws.Range("L27").Value = 7500
    $ws.Range("M27").Value = 7250
    $ws.Range("O27").Value = 'Región Metropolitana'
    $ws.Range("P27").Value = 290
    # Row 28 <- values from original row 22
    $ws.Range("D28").Value = 44533
    $ws.Range("J28").Value = 180
    $ws.Range("K28").Value = 8000
    $ws.Range("L28").Value = 8500
    $ws.Range("M28").Value = 8222
    $ws.Range("O28").Value = 'Región del Maule'
    $ws.Range("P28").Value = 329
    # Row 29 <- values from original row 19
    $ws.Range("D29").Value = 44454
    $ws.Range("J29").Value = 100
    $ws.Range("K29").Value = 13000
    $ws.Range("L29").Value = 14000
    $ws.Range("M29").Value = 13500
    $ws.Range("O29").Value = 'Provincia del Elquí'
    $ws.Range("P29").Value = 540
    # Row 31 <- values from original row 26
    $ws.Range("D31").Value = 44782
    $ws.Range("J31").Value = 30
    $ws.Range("K31").Value = 15000
    $ws.Range("L31").Value = 15000
    $ws.Range("M31").Value = 15000
    $ws.Range("O31").Value = 'Región de Coquimbo'
    $ws.Range("P31").Value = 600

